$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.078.04"
$ws.Range("E2").Value = "  +4.30%  "

# Row 3
$ws.Range("D3").Value = "1.691.87"
$ws.Range("E3").Value = "  +3.33%  "

# Row 4
$ws.Range("D4").Value = "'0.9955"
$ws.Range("E4").Value = "  -0.37%  "

# Row 5
$ws.Range("D5").Value = "'239.68"
$ws.Range("E5").Value = "  +3.45%  "

# Row 6
$ws.Range("D6").Value = "'0.9969"
$ws.Range("E6").Value = "  -0.34%  "

# Row 7
$ws.Range("D7").Value = "'0.4666"
$ws.Range("E7").Value = "  -1.06%  "

# Row 8
$ws.Range("D8").Value = "'0.2624"
$ws.Range("E8").Value = "  +3.02%  "

# Row 9
$ws.Range("D9").Value = "'0.06178"
$ws.Range("E9").Value = "  +1.98%  "

# Row 10
$ws.Range("D10").Value = "1.678.28"
$ws.Range("E10").Value = "  +2.25%  "

# Row 11
$ws.Range("D11").Value = "'0.07034"
$ws.Range("E11").Value = "  +0.21%  "

# Row 12
$ws.Range("D12").Value = "'15.20"
$ws.Range("E12").Value = "  +6.27%  "

# Row 13
$ws.Range("D13").Value = "'4.412"
$ws.Range("E13").Value = "  +3.04%  "

# Row 14
$ws.Range("D14").Value = "'0.5858"
$ws.Range("E14").Value = "  +3.19%  "

# Row 15
$ws.Range("D15").Value = "'75.94"
$ws.Range("E15").Value = "  +3.80%  "

# Row 16
$ws.Range("D16").Value = "'0.9991"
$ws.Range("E16").Value = "  -0.14%  "

# Row 17
$ws.Range("D17").Value = "'0.9963"
$ws.Range("E17").Value = "  -0.38%  "

# Row 18
$ws.Range("D18").Value = "26.032.14"
$ws.Range("E18").Value = "  +4.17%  "

# Row 19
$ws.Range("D19").Value = "'0.000006753"
$ws.Range("E19").Value = "  +3.23%  "

# Row 20
$ws.Range("D20").Value = "'11.51"
$ws.Range("E20").Value = "  +3.01%  "

# Row 21
$ws.Range("D21").Value = "1.894.61"
$ws.Range("E21").Value = "  +2.33%  "

# Row 22
$ws.Range("D22").Value = "'4.522"
$ws.Range("E22").Value = "  +6.42%  "

# Row 23
$ws.Range("D23").Value = "'8.742"
$ws.Range("E23").Value = "  +3.45%  "

# Row 24
$ws.Range("D24").Value = "'5.276"
$ws.Range("E24").Value = "  +1.80%  "

# Row 25
$ws.Range("D25").Value = "'134.66"
$ws.Range("E25").Value = "  +1.79%  "

# Row 26
$ws.Range("D26").Value = "'15.07"
$ws.Range("E26").Value = "  +1.95%  "

# Row 27
$ws.Range("D27").Value = "'1.439"
$ws.Range("E27").Value = "  +5.85%  "

# Row 28
$ws.Range("D28").Value = "'1.732"
$ws.Range("E28").Value = "  +6.66%  "

# Row 29
$ws.Range("D29").Value = "'105.42"
$ws.Range("E29").Value = "  +2.14%  "

# Row 30
$ws.Range("D30").Value = "'3.977"
$ws.Range("E30").Value = "  +2.86%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'3.658"
$ws.Range("E31").Value = "  +4.40%  "

# Row 32
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.07751"
$ws.Range("E32").Value = "  +2.69%  "

# Row 33
$ws.Range("D33").Value = "'0.04377"
$ws.Range("E33").Value = "  +3.82%  "

# Row 34
$ws.Range("D34").Value = "'2.613"
$ws.Range("E34").Value = "  +1.71%  "

# Row 35
$ws.Range("D35").Value = "'0.9597"
$ws.Range("E35").Value = "  +3.64%  "

# Row 36
$ws.Range("D36").Value = "'0.6115"
$ws.Range("E36").Value = "  +4.16%  "

# Row 37
$ws.Range("D37").Value = "'0.9318"
$ws.Range("E37").Value = "  +4.82%  "

# Row 38
$ws.Range("E38").Value = "  +13.59%  "

# Row 39
$ws.Range("D39").Value = "'2.391"
$ws.Range("E39").Value = "  -6.75%  "

# Row 40
$ws.Range("D40").Value = "'0.9972"
$ws.Range("E40").Value = "  -0.28%  "

# Row 41
$ws.Range("D41").Value = "'1.893"
$ws.Range("E41").Value = "  +7.97%  "

# Row 42
$ws.Range("D42").Value = "'0.01465"
$ws.Range("E42").Value = "  -0.83%  "

# Row 43
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.3756"
$ws.Range("E43").Value = "  +2.95%  "

# Row 44
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.072"
$ws.Range("E44").Value = "  +9.68%  "

# Row 45
$ws.Range("D45").Value = "'0.1131"
$ws.Range("E45").Value = "  +3.83%  "

# Row 46
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "'6.209"
$ws.Range("E46").Value = "  +2.86%  "

# Row 47
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.05310"
$ws.Range("E47").Value = "  +2.33%  "

# Row 48
$ws.Range("E48").Value = "  +8.38%  "

# Row 49
$ws.Range("D49").Value = "'7.653"
$ws.Range("E49").Value = "  +8.96%  "

# Row 50
$ws.Range("D50").Value = "'1.216"
$ws.Range("E50").Value = "  +2.94%  "

# Row 51
$ws.Range("D51").Value = "'0.9990"
$ws.Range("E51").Value = "  -0.24%  "
